$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" column (E) for rows 16-26 currently lists periods in
# descending order (1808 down to 1710). Update the database so the new
# account-statement rows are listed in ascending order (1710 up to 1808),
# matching the refreshed EC database / part 1 of new statements.
$periodos = @("1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
